$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Uvar")
$ws.Activate()

# Update risk scenario values in column B (preliminary risk scenarios)
$ws.Range("B3").Value = -0.10513400000000001
$ws.Range("B4").Value = 0.078451000000000007
$ws.Range("B5").Value = 0.21173500000000001
$ws.Range("B6").Value = 0.27
$ws.Range("B7").Value = 0.3
$ws.Range("B8").Value = 0.3
$ws.Range("B9").Value = 0.3
$ws.Range("B10").Value = 0.3

# Update the selected cell to reflect the last-edited location
$ws.Range("H8").Select()
